# quarterly.xlsx "Overview" update
#
# The source workbook tracks 10 rolling fiscal quarters (columns E:N) for a
# handful of G&A expense line items. This edit rolls the window forward by
# one quarter: the oldest quarter column ("فصل دوم منتهی به 1399/06") is
# dropped and a new quarter ("فصل چهارم منتهی به 1401/12") is appended,
# with every existing quarter's figures shifting one column to the left
# (E<-F, F<-G, ... M<-N) and a freshly reported value landing in N.
#
# This mirrors exactly what a user would do by hand in Excel: retype the
# header labels in row 8 / row 24 and retype the figures in each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# --- Rolling quarter header labels (row 8 and row 24) ---
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

foreach ($headerRow in @(8, 24)) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $headerRow).Value = $quarterLabels[$i]
    }
}

# --- Quarterly figures, shifted one column left with a new value in N ---
$dataRows = @{
    10 = @(88315, 94784, 101804, 248322, 115234, 80063, 305062, 199376, 179079, 51091)
    11 = @(177997, 186533, 280148, 122253, 160239, 274446, 284743, 251460, 301995, 434005)
    12 = @(13043, 16025, 30931, 30503, 46345, 11938, 96936, 22692, 61615, 70197)
    16 = @(8733, 9291, 8862, 8964, 9604, 10091, 9657, 9741, 10231, 16181)
    17 = @(99309, 160871, 138562, 112278, 141749, 188947, 235772, 201731, 263037, 352680)
    19 = @(21454, 17235, 17774, 20969, 45571, 53548, 40996, 81949, 136275, 148670)
    20 = @(408851, 484739, 578081, 543289, 518742, 619033, 973166, 766949, 952232, 1072824)
    26 = @(192, 193, 189, 191, 191, 190, 189, 188, 188, 186)
    27 = @(312, 311, 309, 306, 306, 303, 299, 299, 298, 297)
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# --- Row heights: the rows that already carried an explicit height shrank
#     slightly (15 -> 14.4 default-row-height scale, i.e. x0.96) ---
$rowHeights = @{
    2  = 15.6
    5  = 40.8
    6  = 40.8
    8  = 31.2
    24 = 31.2
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}
